$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn thu nợ")

# Row 3 used to be the "Tổng" (summary) row for the single data row (row 2).
# It becomes a new data row ("TN") with its own values.
$ws.Cells.Item(3, 1).Value = "TN"
$ws.Cells.Item(3, 2).Value = 160
$ws.Cells.Item(3, 3).Value = 1000000

# Force column D to text so the date-looking string isn't auto-converted
# into a date serial number, then clear the leftover number format style.
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "07-20-2024"
$ws.Cells.Item(3, 4).Style = "Normal"

$ws.Cells.Item(3, 5).Value = "SÓC TRĂNG"
$ws.Cells.Item(3, 6).Value = "HD-LUXURY-488"
$ws.Cells.Item(3, 7).Value = "Nâng mũi"
$ws.Cells.Item(3, 8).Value = "khanh ktv cũ"
$ws.Cells.Item(3, 9).Value = "Khách cửa hàng"
$ws.Cells.Item(3, 10).ClearContents()
$ws.Cells.Item(3, 11).Value = 7900000
$ws.Cells.Item(3, 12).ClearContents()
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 14).Value = 7900000
$ws.Cells.Item(3, 15).Value = 1000000
$ws.Cells.Item(3, 16).Value = "Bác Sĩ Ngoài"
$ws.Cells.Item(3, 17).ClearContents()
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0.1
$ws.Cells.Item(3, 23).Value = 100000
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(3, 25).Value = 0

# Row 4 is the new "Tổng" (summary) row totalling rows 2 and 3.
$ws.Cells.Item(4, 1).Value = "Tổng"
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = 3000000
$ws.Cells.Item(4, 4).ClearContents()
$ws.Cells.Item(4, 5).ClearContents()
$ws.Cells.Item(4, 6).ClearContents()
$ws.Cells.Item(4, 7).ClearContents()
$ws.Cells.Item(4, 8).ClearContents()
$ws.Cells.Item(4, 9).ClearContents()
$ws.Cells.Item(4, 10).ClearContents()
$ws.Cells.Item(4, 11).Value = 14900000
$ws.Cells.Item(4, 12).ClearContents()
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 14900000
$ws.Cells.Item(4, 15).Value = 8000000
$ws.Cells.Item(4, 16).ClearContents()
$ws.Cells.Item(4, 17).ClearContents()
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0.2
$ws.Cells.Item(4, 23).Value = 300000
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0
